$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns for data rows (2-51) are text-formatted so that
# numeric-looking strings (e.g. "1.001", "0.000007579") are preserved exactly
# as text rather than being converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.542.09'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '1.918.61'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '245.39'
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.4789'
$ws.Range("E7").Value = '  +1.62%  '
$ws.Range("D8").Value = '0.2902'
$ws.Range("E8").Value = '  +0.60%  '
$ws.Range("D9").Value = '0.06722'
$ws.Range("E9").Value = '  -0.98%  '
$ws.Range("D10").Value = '110.72'
$ws.Range("E10").Value = '  +3.81%  '
$ws.Range("D11").Value = '19.02'
$ws.Range("E11").Value = '  +3.36%  '
$ws.Range("D12").Value = '1.915.23'
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").Value = '0.07573'
$ws.Range("E13").Value = '  -2.47%  '
$ws.Range("D14").Value = '5.288'
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").Value = '0.6706'
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").Value = '298.90'
$ws.Range("E16").Value = '  +1.79%  '
$ws.Range("D17").Value = '30.545.95'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '5.627'
$ws.Range("E18").Value = '  +4.93%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '0.9998'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '12.96'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = '0.000007579'
$ws.Range("E21").Value = '  -0.39%  '
$ws.Range("D22").Value = '2.166.46'
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '6.520'
$ws.Range("E24").Value = '  +4.52%  '
$ws.Range("D25").Value = '9.483'
$ws.Range("E25").Value = '  +1.14%  '
$ws.Range("D26").Value = '164.68'
$ws.Range("E26").Value = '  -2.49%  '
$ws.Range("D27").Value = '20.25'
$ws.Range("E27").Value = '  -5.22%  '
$ws.Range("D28").Value = '2.108'
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("D29").Value = '0.1072'
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").Value = '1.399'
$ws.Range("E30").Value = '  +2.37%  '
$ws.Range("D31").Value = '4.154'
$ws.Range("E31").Value = '  -1.10%  '
$ws.Range("D32").Value = '4.048'
$ws.Range("E32").Value = '  +0.95%  '
$ws.Range("D33").Value = '0.05011'
$ws.Range("E33").Value = '  -0.79%  '
$ws.Range("D34").Value = '0.7396'
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("D35").Value = '1.140'
$ws.Range("E35").Value = '  -1.45%  '
$ws.Range("D36").Value = '0.9997'
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").Value = '2.734'
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").Value = '0.02026'
$ws.Range("E38").Value = '  -3.87%  '
$ws.Range("D39").Value = '2.685'
$ws.Range("D40").Value = '111.18'
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").Value = '2.027'
$ws.Range("E41").Value = '  -2.60%  '
$ws.Range("D42").Value = '0.4476'
$ws.Range("E42").Value = '  +4.32%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.8646'
$ws.Range("E43").Value = '  -1.70%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '71.58'
$ws.Range("E44").Value = '  +5.39%  '
$ws.Range("D45").Value = '5.878'
$ws.Range("E45").Value = '  -1.01%  '
$ws.Range("D46").Value = '0.9992'
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").Value = '49.53'
$ws.Range("E47").Value = '  -0.71%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.265'
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("D49").Value = '9.283'
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("D50").Value = '0.2580'
$ws.Range("E50").Value = '  +4.54%  '
$ws.Range("D51").Value = '0.1238'
$ws.Range("E51").Value = '  +1.30%  '
